$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 113.333336
$ws.Range("I9").Value = 137.25
$ws.Range("J9").Value = 65.5
$ws.Range("K9").Value = 137.25
$ws.Range("L9").Value = 65.5
$ws.Range("M9").Value = 31.75
$ws.Range("N9").Value = -403.5
$ws.Range("H17").Value = 494.9091
$ws.Range("J17").Value = 494.9091
$ws.Range("L17").Value = 1484.7273
$ws.Range("N17").Value = -1820.7273
$ws.Range("H19").Value = 591.05884
$ws.Range("I19").Value = 585
$ws.Range("J19").Value = 595.3
$ws.Range("K19").Value = 585
$ws.Range("L19").Value = 595.3
$ws.Range("M19").Value = -410
$ws.Range("N19").Value = -945.3
$ws.Range("H86").Value = 2120.8333
$ws.Range("I86").Value = 1730.5625
$ws.Range("J86").Value = 2901.375
$ws.Range("K86").Value = 1730.5625
$ws.Range("L86").Value = 2901.375
$ws.Range("M86").Value = -607.5625
$ws.Range("N86").Value = -5147.375
$ws.Range("H89").Value = 2120.8333
$ws.Range("I89").Value = 1730.5625
$ws.Range("J89").Value = 2901.375
$ws.Range("K89").Value = 8652.8125
$ws.Range("L89").Value = 14506.875
$ws.Range("M89").Value = -3036.8125
$ws.Range("N89").Value = -25738.875
$ws.Range("H112").Value = 2111.077
$ws.Range("I112").Value = 539.25
$ws.Range("J112").Value = 2396.8635
$ws.Range("K112").Value = 1617.75
$ws.Range("L112").Value = 7190.5905
$ws.Range("M112").Value = -509.75
$ws.Range("N112").Value = -9406.5905
$ws.Range("H127").Value = 1291.6
$ws.Range("I127").Value = 846.75
$ws.Range("J127").Value = 1453.3636
$ws.Range("K127").Value = 2540.25
$ws.Range("L127").Value = 4360.0908
$ws.Range("M127").Value = 2419.75
$ws.Range("N127").Value = -14280.0908
$ws.Range("H129").Value = 1146.5283
$ws.Range("I129").Value = 463.4
$ws.Range("K129").Value = 1390.2
$ws.Range("M129").Value = 3609.8
$ws.Range("H138").Value = 4162.7856
$ws.Range("I138").Value = 1929.3334
$ws.Range("J138").Value = 4771.909
$ws.Range("K138").Value = 5788.0002
$ws.Range("L138").Value = 14315.727
$ws.Range("M138").Value = -648.0002000000004
$ws.Range("N138").Value = -24595.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -84
$ws.Range("H28").Value = 500000
$ws.Range("I28").Value = 500000
$ws.Range("K28").Value = 500000
$ws.Range("M28").Value = -499808
$ws.Range("H45").Value = 1166.6666
$ws.Range("I45").Value = 1166.6666
$ws.Range("K45").Value = 1166.6666
$ws.Range("M45").Value = -789.6666
$ws.Range("H74").Value = 1793
$ws.Range("I74").Value = 1940.9
$ws.Range("J74").Value = 1300
$ws.Range("K74").Value = 1940.9
$ws.Range("L74").Value = 1300
$ws.Range("M74").Value = -1066.9
$ws.Range("N74").Value = -3048
$ws.Range("H77").Value = 1793
$ws.Range("I77").Value = 1940.9
$ws.Range("J77").Value = 1300
$ws.Range("K77").Value = 9704.5
$ws.Range("L77").Value = 6500
$ws.Range("M77").Value = -5336.5
$ws.Range("N77").Value = -15236
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H99").Value = 500000
$ws.Range("I99").Value = 500000
$ws.Range("K99").Value = 500000
$ws.Range("M99").Value = -497005
$ws.Range("H132").Value = 3804.739
$ws.Range("I132").Value = 3124.75
$ws.Range("K132").Value = 9374.25
$ws.Range("M132").Value = -6844.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 170335
$ws.Range("I86").Value = 4751.25
$ws.Range("J86").Value = 501502.5
$ws.Range("K86").Value = 4751.25
$ws.Range("L86").Value = 501502.5
$ws.Range("M86").Value = -3628.25
$ws.Range("N86").Value = -503748.5
$ws.Range("H89").Value = 170335
$ws.Range("I89").Value = 4751.25
$ws.Range("J89").Value = 501502.5
$ws.Range("K89").Value = 23756.25
$ws.Range("L89").Value = 2507512.5
$ws.Range("M89").Value = -18140.25
$ws.Range("N89").Value = -2518744.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2291.5
$ws.Range("I31").Value = 1320.3077
$ws.Range("J31").Value = 6500
$ws.Range("K31").Value = 1320.3077
$ws.Range("L31").Value = 6500
$ws.Range("M31").Value = -1025.3077
$ws.Range("N31").Value = -7090
$ws.Range("H34").Value = 2291.5
$ws.Range("I34").Value = 1320.3077
$ws.Range("J34").Value = 6500
$ws.Range("K34").Value = 1320.3077
$ws.Range("L34").Value = 6500
$ws.Range("M34").Value = -1118.3077
$ws.Range("N34").Value = -6904
$ws.Range("H62").Value = 48082.273
$ws.Range("I62").Value = 57900.555
$ws.Range("K62").Value = 57900.555
$ws.Range("M62").Value = -57276.555
$ws.Range("H65").Value = 48082.273
$ws.Range("I65").Value = 57900.555
$ws.Range("K65").Value = 289502.775
$ws.Range("M65").Value = -286382.775
$ws.Range("H86").Value = 2812.6365
$ws.Range("I86").Value = 1989.5
$ws.Range("J86").Value = 2995.5557
$ws.Range("K86").Value = 1989.5
$ws.Range("L86").Value = 2995.5557
$ws.Range("M86").Value = -866.5
$ws.Range("N86").Value = -5241.5557
$ws.Range("H89").Value = 2812.6365
$ws.Range("I89").Value = 1989.5
$ws.Range("J89").Value = 2995.5557
$ws.Range("K89").Value = 9947.5
$ws.Range("L89").Value = 14977.7785
$ws.Range("M89").Value = -4331.5
$ws.Range("N89").Value = -26209.7785

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 50012
$ws.Range("I70").Value = 50012
$ws.Range("K70").Value = 150036
$ws.Range("M70").Value = -149721
$ws.Range("H73").Value = 50012
$ws.Range("I73").Value = 50012
$ws.Range("K73").Value = 150036
$ws.Range("M73").Value = -148944
$ws.Range("H113").Value = 671.6042
$ws.Range("I113").Value = 598.2414
$ws.Range("J113").Value = 783.5789
$ws.Range("K113").Value = 1794.7242
$ws.Range("L113").Value = 2350.7367
$ws.Range("M113").Value = 375.2757999999999
$ws.Range("N113").Value = -6690.736699999999
$ws.Range("H133").Value = 4533.222
$ws.Range("I133").Value = 1966
$ws.Range("J133").Value = 7742.25
$ws.Range("K133").Value = 5898
$ws.Range("L133").Value = 23226.75
$ws.Range("M133").Value = -838
$ws.Range("N133").Value = -33346.75
$ws.Range("H138").Value = 3042
$ws.Range("I138").Value = 1665
$ws.Range("J138").Value = 3501
$ws.Range("K138").Value = 4995
$ws.Range("L138").Value = 10503
$ws.Range("M138").Value = 145
$ws.Range("N138").Value = -20783

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7878.227
$ws.Range("I70").Value = 8781.4
$ws.Range("J70").Value = 5942.857
$ws.Range("K70").Value = 8781.4
$ws.Range("L70").Value = 5942.857
$ws.Range("M70").Value = -8511.4
$ws.Range("N70").Value = -6482.857
$ws.Range("H73").Value = 7878.227
$ws.Range("I73").Value = 8781.4
$ws.Range("J73").Value = 5942.857
$ws.Range("K73").Value = 8781.4
$ws.Range("L73").Value = 5942.857
$ws.Range("M73").Value = -7845.4
$ws.Range("N73").Value = -7814.857
$ws.Range("H102").Value = 3190.2
$ws.Range("I102").Value = 3218.6667
$ws.Range("J102").Value = 3147.5
$ws.Range("K102").Value = 3218.6667
$ws.Range("L102").Value = 3147.5
$ws.Range("M102").Value = -1596.6667
$ws.Range("N102").Value = -6391.5
$ws.Range("H125").Value = 43333.332
$ws.Range("J125").Value = 43333.332
$ws.Range("L125").Value = 43333.332
$ws.Range("N125").Value = -48253.332
$ws.Range("H135").Value = 47481.43
$ws.Range("J135").Value = 47481.43
$ws.Range("L135").Value = 47481.43
$ws.Range("N135").Value = -57621.43
$ws.Range("H137").Value = 75858.336
$ws.Range("J137").Value = 75858.336
$ws.Range("L137").Value = 75858.336
$ws.Range("N137").Value = -86058.336
$ws.Range("H138").Value = 84625
$ws.Range("J138").Value = 84625
$ws.Range("L138").Value = 84625
$ws.Range("N138").Value = -94905
$ws.Range("H140").Value = 88875.45
$ws.Range("J140").Value = 88875.45
$ws.Range("L140").Value = 88875.45
$ws.Range("N140").Value = -99235.45
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 692.8570999999999
$ws.Range("J46").Value = 670
$ws.Range("L46").Value = 670
$ws.Range("N46").Value = -1046

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 2620.1667
$ws.Range("I132").Value = 1734.4286
$ws.Range("K132").Value = 5203.2858
$ws.Range("M132").Value = -2673.2858
